$d = $word.ActiveDocument

# Locate the paragraph/run text that needs to be replaced.
$seek = $d.Content
$seek.Find.Execute(
    "Replace hard coded ids (such as 1) with fragment of code in all tests",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)

if (-not $seek.Find.Found) {
    throw "Could not find the target text to replace"
}

# Re-materialize the found span as a brand-new Range; InsertXML only performs
# a true replace (rather than an append) on a freshly constructed Range.
$target = $d.Range($seek.Start, $seek.End)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:rStyle w:val="HTMLCode"/><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorBidi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>'

$payload = '<?xml version="1.0" standalone="yes"?>' +
    '<?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
          "<w:document $ns>" +
            '<w:body>' +
              '<w:p>' +
                "<w:r>$rPr" + '<w:t xml:space="preserve">Find out whether </w:t></w:r>' +
                '<w:proofErr w:type="spellStart"/>' +
                "<w:r>$rPr" + '<w:t>password_valid</w:t></w:r>' +
                '<w:proofErr w:type="spellEnd"/>' +
                "<w:r>$rPr" + '<w:t>? can be replaced with regular validation in User moder</w:t></w:r>' +
              '</w:p>' +
            '</w:body>' +
          '</w:document>' +
        '</pkg:xmlData>' +
      '</pkg:part>' +
    '</pkg:package>'

$target.InsertXML($payload)
